$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Merge the three runs that make up the "third part of the staff system"
#    paragraph into a single run (no visible text change, just simplifies the
#    run structure the same way Word does when the paragraph is re-typed).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "The third part of the staff system determines the distribution of attention of staff members. Distribution of water through tankers in the case of a water crisis always takes priority. If after that enough staff members are available for all work, no distribution is made. However, when staff capacity is limited 25% of staff attention is given to maintenance activities, and 75% to refurbishing and planning activities. The refurbishing and planning activities are divided pro ratio.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The third part of the staff system determines the distribution of attention of staff members. Distribution of water through tankers in the case of a water crisis always takes priority. If after that enough staff members are available for all work, no distribution is made. However, when staff capacity is limited 25% of staff attention is given to maintenance activities, and 75% to refurbishing and planning activities. The refurbishing and planning activities are divided pro ratio.",
    2
) | Out-Null

# ---------------------------------------------------------------------------
# 2. Rework the "Effect of maintenance on ageing time" conclusion sentence.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Having a higher maintenance coverage than 8% will not increase the lifespan of infrastructure.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "If all infrastructure is maintained each year, the average infrastructure aging time will increase by 10 years.",
    2
) | Out-Null

# ---------------------------------------------------------------------------
# 3. Merge the "In the " / "base case scenario the model shows" runs (this is
#    also where the stray _GoBack bookmark used to sit) into one run.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "In the base case scenario the model shows",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "In the base case scenario the model shows",
    2
) | Out-Null

# ---------------------------------------------------------------------------
# 4. Re-seat the (last-edit) _GoBack bookmark onto the single space between
#    "...the background reading" and "it was described..." in the leakage
#    paragraph, matching where the author's cursor ended up after editing.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("the background reading it was described that there is an average leakage") | Out-Null
$spaceStart = $rng.Start + 22
$spaceEnd = $spaceStart + 1
$bm = $d.Range($spaceStart, $spaceEnd)
$d.Bookmarks.Add("_GoBack", $bm) | Out-Null
